# Apply the edits described by the commit:
# "new .ttl from Google sheet has been generated"
#
# 1. Update the last-updated timestamp (row 17, col B).
# 2. Tweak the English definition of "other" (row 25, col E):
#    "data" -> "resources".
# 3. Remove the "other type of biomaterial" term (row 31) entirely —
#    all following rows (32-93) shift up by one, and the sheet's used
#    range shrinks from AF93 to AF92.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "2021-08-13T22:00:00+00:00"
$ws.Range("E25").Value = "Other forms of resources that are not easily categorized or defined"

$ws.Rows.Item(31).Delete()
